$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header text in E5 so it references the "$ CHG " (trailing space) shared
# string instead of the duplicate "$ CHG" string - this also causes the now-unused
# duplicate shared string to be dropped when the workbook is saved.
$ws.Range("E5").Value = "`$ CHG "

# Update the East row (row 6) figures and turn the $ CHG / % CHG columns into
# formulas driven off the other cells in the row.
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 400
$ws.Range("D6").Value = 300
$ws.Range("E6").Formula = "=B6-C6"
$ws.Range("F6").Formula = "=E6/C6"

# Remove the now unneeded extra data rows (7, 8 and 9).
$ws.Rows("7:9").Delete()
